$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E / B cleanup around rows 10-14 -----------------------------
# "?Wong Ming" (row 13) is resolved into a real entry "Wong Ming" (added
# below, row 19). Row 14's "?Wong Kei" shifts up to fill row 13, leaving
# row 14 empty. "?#Denvendra"/row11-E, "?Casey"/row12-E, "Ray"/row13-E are
# removed from this block (their "reconfirmed" values reappear lower down).
$ws.Range("E10").Value = "?Casey"
$ws.Range("E11").Value = $null
$ws.Range("E12").Value = $null
$ws.Range("B13").Value = "?Wong Kei"
$ws.Range("E13").Value = $null
$ws.Range("B14").Value = $null

# --- Rows 19-25: column D shifts up one slot, column B gains two new ----
# entries ("Wong Ming" confirmed, "Guanglei" moved from D), and column E
# gains the reconfirmed names pulled up from the block above.
$ws.Range("B19").Value = "Wong Ming"
$ws.Range("D19").Value = "Tim"

$ws.Range("B20").Value = "Guanglei"
$ws.Range("D20").Value = "Patrick"

$ws.Range("D21").Value = "Ah Chicken"

$ws.Range("D22").Value = "Bean Man"
$ws.Range("E22").Value = "Ray"

$ws.Range("D23").Value = "Foo Kwai"
$ws.Range("E23").Value = "Shirley"

$ws.Range("D24").Value = "Ellen "
$ws.Range("E24").Value = "#Denvendra"

$ws.Range("D25").Value = $null

# --- Restore the selection recorded in the saved workbook ---------------
$ws.Range("E8").Select()
